$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.758.91"
$ws.Range("E2").Value = "  +3.55%  "

$ws.Range("D3").Value = "3.691.07"
$ws.Range("E3").Value = "  +8.54%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "589.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.52%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "180.08"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.69%  "

$ws.Range("D7").Value = "3.683.54"
$ws.Range("E7").Value = "  +8.56%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.622"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.04%  "

$ws.Range("E9").Value = "  +0.02%  "

$ws.Range("E10").Value = "  +1.68%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.615"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.16%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "50.05"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.48%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000287"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.67%  "

$ws.Range("D14").Value = "4.286.16"
$ws.Range("E14").Value = "  +8.54%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "686.04"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.80%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "9.01"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.73%  "

$ws.Range("D17").Value = "3.696.78"
$ws.Range("E17").Value = "  +8.59%  "

$ws.Range("D18").Value = "71.816.81"
$ws.Range("E18").Value = "  +3.46%  "

$ws.Range("E19").Value = "  +2.26%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.13"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.42%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.69"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.45%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.943"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.74%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.32"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +17.73%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "17.89"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.97%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "104.15"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.23%  "

$ws.Range("E26").Value = "  +3.79%  "

$ws.Range("E27").Value = "  +5.78%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.17"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.27%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "35.57"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.27%  "

$ws.Range("E30").Value = "  +5.71%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.35"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.82%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.25"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +12.70%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "570.17"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.48%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.33"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.75%  "

$ws.Range("E35").Value = "  +4.11%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "59.58"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.60%  "

$ws.Range("D37").Value = "3.809.76"
$ws.Range("E37").Value = "  +5.59%  "

$ws.Range("E38").Value = "  -0.02%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.147"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.87%  "

$ws.Range("D40").Value = "0.0₃0783"
$ws.Range("E40").Value = "  +5.09%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "35.49"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.73%  "

$ws.Range("E42").Value = "  +5.77%  "

$ws.Range("E43").Value = "  +9.22%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.80"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.68%  "

$ws.Range("E45").Value = "  +5.34%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.89"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +8.00%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.37"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.66%  "

$ws.Range("E48").Value = "  +4.25%  "

$ws.Range("E49").Value = "  +2.43%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.998"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.23%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "135.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.78%  "
